# Auto-committed on 2023/09/15 週五 17:07:32.90
#
# Adds a new "更新BY / L7205-五類資產分類上傳轉檔作業" tracking column to the
# DBD layout sheet, switches the two "last updated" timestamp fields from
# DATE to TIMESTAMP, documents the new jobTxSeq stored-proc parameter on the
# SP sheet, and records the new findYearMonthAll lookup on the DBS sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "DBD" - add the new "更新BY..." column (H) and retype the two DATE
# fields (CreateDate / LastUpdate) to TIMESTAMP.
# ---------------------------------------------------------------------------
$dbd = $wb.Worksheets.Item("DBD")
$dbd.Activate()

$dbd.Columns.Item(8).ColumnWidth = 41.109375
$dbd.Range("H1:H4").WrapText = $true
$dbd.Cells.Item(1, 8).Value = "更新BY" + [char]10 + "L7205-五類資產分類上傳轉檔作業"
$dbd.Rows.Item(1).RowHeight = 32.4

$dbd.Cells.Item(15, 4).Value = "TIMESTAMP"
$dbd.Cells.Item(17, 4).Value = "TIMESTAMP"

$dbd.Range("D17").Select()

# ---------------------------------------------------------------------------
# Sheet "DBS" - document the new findYearMonthAll lookup function.
# ---------------------------------------------------------------------------
$dbs = $wb.Worksheets.Item("DBS")
$dbs.Activate()

$dbs.Cells.Item(2, 1).Value = "findYearMonthAll"
$dbs.Cells.Item(2, 2).Value = "YearMonth = "

$dbs.Range("A2").Activate()
$dbs.Range("A2:B2").Select()

# ---------------------------------------------------------------------------
# Sheet "SP" - the insert stored procedure now also takes a jobTxSeq param.
# ---------------------------------------------------------------------------
$sp = $wb.Worksheets.Item("SP")
$sp.Activate()

$sp.Cells.Item(2, 2).Value = "int tbsdyf,  String empNo,  String jobTxSeq"

$sp.Range("C5").Select()

# ---------------------------------------------------------------------------
# Leave the DBD sheet active/selected, matching the saved workbook state.
# ---------------------------------------------------------------------------
$dbd.Activate()
